$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '30.461.81'
Set-TextValue $ws.Range('E2') '  -1.14%  '

Set-TextValue $ws.Range('D3') '2.107.52'
Set-TextValue $ws.Range('E3') '  -0.37%  '

Set-TextValue $ws.Range('D4') '1.003'
Set-TextValue $ws.Range('E4') '  +0.17%  '

Set-TextValue $ws.Range('D5') '333.22'
Set-TextValue $ws.Range('E5') '  -0.05%  '

Set-TextValue $ws.Range('E6') '  +0.04%  '

Set-TextValue $ws.Range('D7') '0.5230'
Set-TextValue $ws.Range('E7') '  -1.64%  '

Set-TextValue $ws.Range('D8') '0.4487'
Set-TextValue $ws.Range('E8') '  +2.53%  '

Set-TextValue $ws.Range('D9') '53.61'
Set-TextValue $ws.Range('E9') '  +15.92%  '

Set-TextValue $ws.Range('D10') '0.08981'
Set-TextValue $ws.Range('E10') '  -0.31%  '

Set-TextValue $ws.Range('D11') '1.161'
Set-TextValue $ws.Range('E11') '  -1.58%  '

Set-TextValue $ws.Range('D12') '24.46'
Set-TextValue $ws.Range('E12') '  -2.54%  '

Set-TextValue $ws.Range('D13') '2.105.83'
Set-TextValue $ws.Range('E13') '  -0.03%  '

Set-TextValue $ws.Range('D14') '6.782'
Set-TextValue $ws.Range('E14') '  +0.39%  '

Set-TextValue $ws.Range('D15') '7.786'

Set-TextValue $ws.Range('D16') '96.53'
Set-TextValue $ws.Range('E16') '  -0.99%  '

Set-TextValue $ws.Range('E17') '  +0.07%  '

Set-TextValue $ws.Range('D18') '0.00001125'
Set-TextValue $ws.Range('E18') '  -0.20%  '

Set-TextValue $ws.Range('D19') '0.06613'
Set-TextValue $ws.Range('E19') '  -0.69%  '

Set-TextValue $ws.Range('D20') '19.32'
Set-TextValue $ws.Range('E20') '  +1.06%  '

Set-TextValue $ws.Range('D21') '1.000'
Set-TextValue $ws.Range('E21') '  +0.03%  '

Set-TextValue $ws.Range('D22') '6.304'
Set-TextValue $ws.Range('E22') '  -0.63%  '

Set-TextValue $ws.Range('D23') '30.519.53'
Set-TextValue $ws.Range('E23') '  -1.12%  '

Set-TextValue $ws.Range('D24') '12.33'
Set-TextValue $ws.Range('E24') '  -0.29%  '

Set-TextValue $ws.Range('D25') '2.341'
Set-TextValue $ws.Range('E25') '  +3.15%  '

Set-TextValue $ws.Range('D26') '2.353.43'
Set-TextValue $ws.Range('E26') '  -0.07%  '

Set-TextValue $ws.Range('D27') '22.40'
Set-TextValue $ws.Range('E27') '  -1.70%  '

Set-TextValue $ws.Range('D28') '2.588'
Set-TextValue $ws.Range('E28') '  +0.36%  '

Set-TextValue $ws.Range('D29') '163.65'
Set-TextValue $ws.Range('E29') '  +0.44%  '

Set-TextValue $ws.Range('D30') '133.12'
Set-TextValue $ws.Range('E30') '  -0.20%  '

Set-TextValue $ws.Range('D31') '1.202'
Set-TextValue $ws.Range('E31') '  +2.50%  '

Set-TextValue $ws.Range('D32') '0.1073'
Set-TextValue $ws.Range('E32') '  -0.60%  '

Set-TextValue $ws.Range('D33') '1.668'
Set-TextValue $ws.Range('E33') '  +7.73%  '

Set-TextValue $ws.Range('D34') '6.164'
Set-TextValue $ws.Range('E34') '  -0.97%  '

Set-TextValue $ws.Range('D35') '3.934'
Set-TextValue $ws.Range('E35') '  -2.07%  '

Set-TextValue $ws.Range('D36') '10.52'
Set-TextValue $ws.Range('E36') '  +10.46%  '

Set-TextValue $ws.Range('D37') '0.02573'
Set-TextValue $ws.Range('E37') '  -1.66%  '

Set-TextValue $ws.Range('D38') '0.06796'
Set-TextValue $ws.Range('E38') '  +0.88%  '

Set-TextValue $ws.Range('D39') '12.80'
Set-TextValue $ws.Range('E39') '  -0.55%  '

Set-TextValue $ws.Range('D40') '5.495'
Set-TextValue $ws.Range('E40') '  -0.74%  '

Set-TextValue $ws.Range('D41') '0.2281'
Set-TextValue $ws.Range('E41') '  -0.18%  '

Set-TextValue $ws.Range('D42') '0.6932'
Set-TextValue $ws.Range('E42') '  +1.24%  '

Set-TextValue $ws.Range('D43') '1.255'
Set-TextValue $ws.Range('E43') '  +0.23%  '

Set-TextValue $ws.Range('D44') '2.348'
Set-TextValue $ws.Range('E44') '  +5.28%  '

Set-TextValue $ws.Range('D45') '1.000'
Set-TextValue $ws.Range('E45') '  +0.06%  '

Set-TextValue $ws.Range('B46') 'EnergySwap'
Set-TextValue $ws.Range('C46') 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range('D46') '14.08'
Set-TextValue $ws.Range('E46') '  +0.04%  '

Set-TextValue $ws.Range('B47') 'Decentraland'
Set-TextValue $ws.Range('C47') 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextValue $ws.Range('D47') '0.6385'
Set-TextValue $ws.Range('E47') '  -0.83%  '

Set-TextValue $ws.Range('D48') '3.643'
Set-TextValue $ws.Range('E48') '  -0.48%  '

Set-TextValue $ws.Range('D49') '1.247'
Set-TextValue $ws.Range('E49') '  -2.22%  '

Set-TextValue $ws.Range('D50') '1.221'
Set-TextValue $ws.Range('E50') '  +5.71%  '

Set-TextValue $ws.Range('D51') '83.32'
Set-TextValue $ws.Range('E51') '  +0.67%  '
